$d = $word.ActiveDocument

function Get-ExactParagraph($doc, $exactText) {
    $want = $exactText + "`r"
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        $p = $doc.Paragraphs($i)
        if ($p.Range.Text -eq $want) {
            return $p
        }
    }
    return $null
}

function Replace-InParagraph($doc, $exactOldText, $newText) {
    $p = Get-ExactParagraph $doc $exactOldText
    if ($p -ne $null) {
        $p.Range.Find.Execute($exactOldText, $true, $false, $false, $false, $false,
            $true, 1, $false, $newText, 2) | Out-Null
    }
}

# ---------------------------------------------------------------------------
# 1. Title: "Play Stunt Stars Free: A Thrilling Action-Packed Slot Game"
#    -> "Play Stunt Stars for Free - Action-Packed Online Slot Game"
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "Play Stunt Stars Free: A Thrilling Action-Packed Slot Game",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Play Stunt Stars for Free - Action-Packed Online Slot Game", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2. Build the new bold "Play Stunt Stars for Free - Action-Packed Online
#    Slot Game" paragraph that now appears right before the closing
#    image-prompt paragraph. We reuse the formatting of the bold
#    "Meta description" run (still present at this point) so the new
#    paragraph picks up matching bold styling, then retarget its text.
#    (Paragraph.Next / .Previous are not reliable in this runtime, so we
#    navigate strictly via Paragraphs(index) / Paragraphs.Count.)
# ---------------------------------------------------------------------------
$metaIndex = 2
$metaPara  = $d.Paragraphs($metaIndex)
$afterMeta = $d.Paragraphs($metaIndex + 1)
$metaRangeFull = $d.Range($metaPara.Range.Start, $afterMeta.Range.Start)
$metaFormatted = $metaRangeFull.FormattedText

$lastIndex = $d.Paragraphs.Count
$lastPara = $d.Paragraphs($lastIndex)
$insertPos = $lastPara.Range.Start
$insertSpot = $d.Range($insertPos, $insertPos)
$insertSpot.FormattedText = $metaFormatted

# The copy added one paragraph; the new paragraph now sits right before the
# (shifted) last paragraph.
$newIndex = $lastIndex
$newPara = $d.Paragraphs($newIndex)
$newScope = $newPara.Range
$newScope.Find.Execute(
    "Meta description", $true, $false, $false, $false, $false, $true, 1, $false,
    "Play Stunt Stars for Free - Action-Packed Online Slot Game", 2) | Out-Null

$newPara2 = $d.Paragraphs($newIndex)
$newTitleLen = "Play Stunt Stars for Free - Action-Packed Online Slot Game".Length
$remainderStart = $newPara2.Range.Start + $newTitleLen
$remainderEnd = $newPara2.Range.End - 1
if ($remainderEnd -gt $remainderStart) {
    $remRange = $d.Range($remainderStart, $remainderEnd)
    $remRange.Delete()
}

# ---------------------------------------------------------------------------
# 3. Remove the original "Meta description: ..." paragraph entirely.
# ---------------------------------------------------------------------------
$d.Paragraphs(2).Range.Delete()

# ---------------------------------------------------------------------------
# 4. "What we like" bullets
# ---------------------------------------------------------------------------
Replace-InParagraph $d "Exciting and immersive gameplay" "Action-packed gameplay with thrilling features"
Replace-InParagraph $d "High-quality symbols and graphics" "High-quality symbols and immersive music"
Replace-InParagraph $d "243 ways to win" "Chance to secure large payouts"
Replace-InParagraph $d "Skydive Bonus and other special symbols" "Demo version available for free play and exploration"

# ---------------------------------------------------------------------------
# 5. "What we don't like" bullets
# ---------------------------------------------------------------------------
Replace-InParagraph $d "High volatility" "Limited bonus features"
Replace-InParagraph $d "No progressive jackpot" "High volatility may not suit all players"

# ---------------------------------------------------------------------------
# 6. Closing italic paragraph text
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "Please design a feature image for " + [char]34 + "Stunt Stars" + [char]34 + " that includes the following: - A cartoon-style depiction of a happy Maya warrior with glasses - The image should be eye-catching and exciting, reflecting the action-packed nature of the game - The warrior should be depicted as if they are in the midst of a high-flying stunt, perhaps leaping through the air with fire and explosions in the background - The image should feature bright colors and high contrast to grab the attention of potential players.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Experience the thrills of Stunt Stars with its action-packed gameplay and high-quality symbols. Play for free!", 2) | Out-Null

Write-Output "Done"
